$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.753.62'
$ws.Range("D3").Value = '3.369.12'
$ws.Range("E3").Value = '  -0.76%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'567.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.82%  '
$ws.Range("D6").Value = "'135.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.41%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.368.92'
$ws.Range("E8").Value = '  -0.74%  '
$ws.Range("E9").Value = '  -1.33%  '
$ws.Range("E10").Value = '  +1.31%  '
$ws.Range("E11").Value = '  -3.99%  '
$ws.Range("E12").Value = '  -2.66%  '
$ws.Range("D13").Value = '3.942.97'
$ws.Range("E13").Value = '  -0.84%  '
$ws.Range("E14").Value = '  +0.38%  '
$ws.Range("D15").Value = "'25.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("D16").Value = '3.369.95'
$ws.Range("E16").Value = '  -0.67%  '
$ws.Range("E17").Value = '  -4.45%  '
$ws.Range("D18").Value = '60.861.77'
$ws.Range("E18").Value = '  -1.51%  '
$ws.Range("D19").Value = "'5.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.84%  '
$ws.Range("E20").Value = '  -4.54%  '
$ws.Range("D21").Value = "'9.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.78%  '
$ws.Range("D22").Value = "'370.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.00%  '
$ws.Range("D23").Value = '3.507.41'
$ws.Range("E23").Value = '  -0.55%  '
$ws.Range("D24").Value = "'0.545"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.53%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = "'70.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.64%  '
$ws.Range("D27").Value = "'0.0000121"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.40%  '
$ws.Range("E28").Value = '  +9.44%  '
$ws.Range("E29").Value = '  -4.15%  '
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("D31").Value = "'7.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.84%  '
$ws.Range("D32").Value = "'7.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.14%  '
$ws.Range("E33").Value = '  -2.43%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").Value = "'23.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.01%  '
$ws.Range("D36").Value = "'5.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.69%  '
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").Value = "'6.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.37%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = "'1.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.31%  '
$ws.Range("D39").Value = "'164.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.56%  '
$ws.Range("D40").Value = "'0.0756"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.48%  '
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("E42").Value = '  -1.52%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = "'25.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.88%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = "'41.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.55%  '
$ws.Range("E45").Value = '  -3.20%  '
$ws.Range("E46").Value = '  -2.81%  '
$ws.Range("E47").Value = '  -6.39%  '
$ws.Range("D48").Value = '2.524.99'
$ws.Range("E48").Value = '  +7.52%  '
$ws.Range("D49").Value = "'23.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.37%  '
$ws.Range("E50").Value = '  -1.74%  '
$ws.Range("D51").Value = "'2.39"
$ws.Range("D51").Style = "Normal"
